$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 922.25
$ws.Range("I2").Value = 897.5
$ws.Range("K2").Value = 897.5
$ws.Range("M2").Value = -784.5
$ws.Range("H62").Value = 6024
$ws.Range("I62").Value = 5286
$ws.Range("K62").Value = 5286
$ws.Range("M62").Value = -4662
$ws.Range("H65").Value = 6024
$ws.Range("I65").Value = 5286
$ws.Range("K65").Value = 26430
$ws.Range("M65").Value = -23310
$ws.Range("H115").Value = 1947.3636
$ws.Range("I115").Value = 1947.3636
$ws.Range("K115").Value = 5842.0908
$ws.Range("M115").Value = -4275.0908
$ws.Range("H132").Value = 1362.2
$ws.Range("I132").Value = 1404.4445
$ws.Range("K132").Value = 4213.333500000001
$ws.Range("M132").Value = -1683.333500000001
$ws.Range("H137").Value = 2038.7894
$ws.Range("I137").Value = 2827.5
$ws.Range("K137").Value = 8482.5
$ws.Range("M137").Value = -5932.5
$ws.Range("H138").Value = 1700.4667
$ws.Range("J138").Value = 2699.8333
$ws.Range("L138").Value = 8099.499899999999
$ws.Range("N138").Value = -18379.4999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2627.0967
$ws.Range("I32").Value = 2627.0967
$ws.Range("K32").Value = 2627.0967
$ws.Range("M32").Value = -2340.0967
$ws.Range("H74").Value = 1776.8182
$ws.Range("I74").Value = 1735.3846
$ws.Range("K74").Value = 1735.3846
$ws.Range("M74").Value = -861.3846000000001
$ws.Range("H77").Value = 1776.8182
$ws.Range("I77").Value = 1735.3846
$ws.Range("K77").Value = 8676.923000000001
$ws.Range("M77").Value = -4308.923000000001
$ws.Range("H122").Value = 1777.2222
$ws.Range("I122").Value = 1650.625
$ws.Range("K122").Value = 4951.875
$ws.Range("M122").Value = -2501.875
$ws.Range("H124").Value = 45214.5
$ws.Range("J124").Value = 45214.5
$ws.Range("L124").Value = 45214.5
$ws.Range("N124").Value = -55034.5
$ws.Range("H125").Value = 8571.666999999999
$ws.Range("J125").Value = 8571.666999999999
$ws.Range("L125").Value = 8571.666999999999
$ws.Range("N125").Value = -18411.667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3053.8125
$ws.Range("I86").Value = 2518.5
$ws.Range("J86").Value = 3589.125
$ws.Range("K86").Value = 2518.5
$ws.Range("L86").Value = 3589.125
$ws.Range("M86").Value = -1395.5
$ws.Range("N86").Value = -5835.125
$ws.Range("H89").Value = 3053.8125
$ws.Range("I89").Value = 2518.5
$ws.Range("J89").Value = 3589.125
$ws.Range("K89").Value = 12592.5
$ws.Range("L89").Value = 17945.625
$ws.Range("M89").Value = -6976.5
$ws.Range("N89").Value = -29177.625
$ws.Range("H105").Value = 2066.2856
$ws.Range("I105").Value = 1610.375
$ws.Range("J105").Value = 2674.1667
$ws.Range("K105").Value = 1610.375
$ws.Range("L105").Value = 2674.1667
$ws.Range("M105").Value = 136.625
$ws.Range("N105").Value = -6168.1667
$ws.Range("H107").Value = 2198.8262
$ws.Range("I107").Value = 1598
$ws.Range("K107").Value = 1598
$ws.Range("M107").Value = 322

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 30.857143
$ws.Range("J15").Value = 34.2
$ws.Range("L15").Value = 102.6
$ws.Range("N15").Value = -382.6
$ws.Range("H25").Value = 8888
$ws.Range("I25").Value = 8888
$ws.Range("K25").Value = 26664
$ws.Range("M25").Value = -26495
$ws.Range("H30").Value = 8888
$ws.Range("I30").Value = 8888
$ws.Range("K30").Value = 26664
$ws.Range("M30").Value = -26562
$ws.Range("H60").Value = 37.666668
$ws.Range("I60").Value = 37.666668
$ws.Range("K60").Value = 113.000004
$ws.Range("M60").Value = 137.999996
$ws.Range("H63").Value = 14656
$ws.Range("I63").Value = 14656
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 43968
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -43219
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 14656
$ws.Range("I66").Value = 14656
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 131904
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -128160
$ws.Range("N66").ClearContents()
$ws.Range("H69").Value = 5400
$ws.Range("J69").Value = 6250
$ws.Range("L69").Value = 18750
$ws.Range("N69").Value = -20372
$ws.Range("H72").Value = 5400
$ws.Range("J72").Value = 6250
$ws.Range("L72").Value = 56250
$ws.Range("N72").Value = -64362
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 699.46155
$ws.Range("J113").Value = 703.5833
$ws.Range("L113").Value = 2110.7499
$ws.Range("N113").Value = -6450.7499
$ws.Range("H123").Value = 3626.6667
$ws.Range("J123").Value = 3558.05
$ws.Range("L123").Value = 10674.15
$ws.Range("N123").Value = -15574.15

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9349.4
$ws.Range("I57").Value = 2527.5
$ws.Range("J57").Value = 13897.333
$ws.Range("K57").Value = 2527.5
$ws.Range("L57").Value = 13897.333
$ws.Range("M57").Value = -1707.5
$ws.Range("N57").Value = -15537.333
$ws.Range("H80").Value = 4893.3335
$ws.Range("I80").Value = 1169
$ws.Range("J80").Value = 6755.5
$ws.Range("K80").Value = 1169
$ws.Range("L80").Value = 6755.5
$ws.Range("M80").Value = -171
$ws.Range("N80").Value = -8751.5
$ws.Range("H83").Value = 4893.3335
$ws.Range("I83").Value = 1169
$ws.Range("J83").Value = 6755.5
$ws.Range("K83").Value = 5845
$ws.Range("L83").Value = 33777.5
$ws.Range("M83").Value = -853
$ws.Range("N83").Value = -43761.5
$ws.Range("H106").Value = 60000
$ws.Range("J106").Value = 60000
$ws.Range("L106").Value = 60000
$ws.Range("N106").Value = -62524
$ws.Range("H122").Value = 3286
$ws.Range("I122").Value = 4323
$ws.Range("J122").Value = 1903.3334
$ws.Range("K122").Value = 12969
$ws.Range("L122").Value = 5710.0002
$ws.Range("M122").Value = -10519
$ws.Range("N122").Value = -10610.0002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3306.3157
$ws.Range("I61").Value = 3226.3572
$ws.Range("K61").Value = 3226.3572
$ws.Range("M61").Value = -3024.3572
$ws.Range("H100").Value = 232643.61
$ws.Range("I100").Value = 274397
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 274397
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -273856
$ws.Range("N100").Value = -4082
$ws.Range("H113").Value = 3306.3157
$ws.Range("I113").Value = 3226.3572
$ws.Range("K113").Value = 3226.3572
$ws.Range("M113").Value = -1056.3572
$ws.Range("H122").Value = 3366.5264
$ws.Range("J122").Value = 3665.625
$ws.Range("L122").Value = 10996.875
$ws.Range("N122").Value = -15896.875
$ws.Range("H123").Value = 37999
$ws.Range("J123").Value = 37999
$ws.Range("L123").Value = 37999
$ws.Range("N123").Value = -47799
$ws.Range("H125").Value = 63333
$ws.Range("J125").Value = 59999.5
$ws.Range("L125").Value = 59999.5
$ws.Range("N125").Value = -69839.5
$ws.Range("H141").Value = 89201.75
$ws.Range("J141").Value = 85519
$ws.Range("L141").Value = 85519
$ws.Range("N141").Value = -95879

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 16332.333
$ws.Range("J115").Value = 16332.333
$ws.Range("L115").Value = 16332.333
$ws.Range("N115").Value = -19466.333
$ws.Range("H122").Value = 3257.3333
$ws.Range("I122").Value = 3257.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9771.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7321.999899999999
$ws.Range("N122").ClearContents()
